$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Marking" row values (per-question mark for right/wrong answers)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Correct the "Total" row values that depend on the marking scheme
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "68 / 112"
